# Icam1-Itgax.xlsx: refresh NATMI LR-pair output with recomputed TPM values.
# The underlying analysis now only emits rows for the "Resolving-Mac" target
# cluster (one row per sending cluster), so the previous 8-row table
# (4 senders x 2 target clusters) collapses to a 4-row table
# (4 senders x 1 target cluster), and every numeric column is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the rows for the old "ECs" target cluster (rows 6-9), which no longer
# exist in the refreshed output. This shrinks the used range to A1:T5.
$ws.Rows("6:9").Delete()

# Row 2 (sender: ECs -> target: Resolving-Mac)
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 44.94181766666667
$ws.Range("N2").Value = 134.825453
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1387.392261721076
$ws.Range("R2").Value = 12486.53035548968
$ws.Range("S2").Value = 0.2985789950947061
$ws.Range("T2").Value = 0.2985789950947061

# Row 3 (sender: FAPs -> target: Resolving-Mac)
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 33.793597
$ws.Range("H3").Value = 101.380791
$ws.Range("I3").Value = 0.3268474027571036
$ws.Range("J3").Value = 0.3268474027571037
$ws.Range("M3").Value = 44.94181766666667
$ws.Range("N3").Value = 134.825453
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1518.745674674814
$ws.Range("R3").Value = 13668.71107207332
$ws.Range("S3").Value = 0.3268474027571036
$ws.Range("T3").Value = 0.3268474027571037

# Row 4 (sender: MuSCs -> target: Resolving-Mac)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 2.981185666666667
$ws.Range("H4").Value = 8.943557
$ws.Range("I4").Value = 0.02883365130639111
$ws.Range("J4").Value = 0.02883365130639111
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 44.94181766666667
$ws.Range("N4").Value = 134.825453
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 133.9799026618135
$ws.Range("R4").Value = 1205.819123956321
$ws.Range("S4").Value = 0.02883365130639111
$ws.Range("T4").Value = 0.02883365130639111

# Row 5 (sender: Resolving-Mac -> target: Resolving-Mac)
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 35.74694633333333
$ws.Range("H5").Value = 107.240839
$ws.Range("I5").Value = 0.3457399508417991
$ws.Range("J5").Value = 0.3457399508417991
$ws.Range("M5").Value = 44.94181766666667
$ws.Range("N5").Value = 134.825453
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1606.532744252786
$ws.Range("R5").Value = 14458.79469827507
$ws.Range("S5").Value = 0.3457399508417991
$ws.Range("T5").Value = 0.3457399508417991
